$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Trade #21 closed at 2026-02-17 08:02:41 - unknown UNKNOWN +0.000%
# A new MarketMaking trade (trade #21, 0-indexed A=21) closed as a small
# loss (-$0.11). This ripples into the "Summary" roll-up stats, the
# "Strategy Status" row for MarketMaking, and appends a new trade row to
# both the "All Trades" log and the per-strategy "MarketMaking" log.
# ---------------------------------------------------------------------------

# --- Summary sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 1199.82    # Current Capital
$ws.Range("B4").Value = -0.18      # Total P&L $
$ws.Range("B5").Value = -0.17      # Total P&L %
$ws.Range("B6").Value = 21         # Total Trades
$ws.Range("B8").Value = 9          # Losing Trades
$ws.Range("B9").Value = 28.57      # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---------------------------
$ws = $wb.Worksheets.Item("Strategy Status")
$ws.Range("C4").Value = 99.81999999999999  # Capital
$ws.Range("D4").Value = 21                 # Trades
$ws.Range("E4").Value = -0.18              # P&L $
$ws.Range("F4").Value = -0.18              # P&L %
$ws.Range("G4").Value = 28.57              # Win Rate %

# --- Append the new trade row (row 22) to a trade log sheet ------------
function Add-TradeRow($ws) {
    $ws.Cells.Item(22, 1).Value = 21
    # Force text so "2026-02-17" isn't auto-converted to a date serial,
    # matching every other date cell in the log (stored as plain text).
    $ws.Cells.Item(22, 2).NumberFormat = "@"
    $ws.Cells.Item(22, 2).Value = "2026-02-17"
    $ws.Cells.Item(22, 3).Value = "08:02:35"
    $ws.Cells.Item(22, 4).Value = "MarketMaking"
    $ws.Cells.Item(22, 5).Value = "DOWN"
    $ws.Cells.Item(22, 6).Value = 0.72
    $ws.Cells.Item(22, 7).Value = 0.61
    $ws.Cells.Item(22, 8).Value = "CLOSED"
    $ws.Cells.Item(22, 9).Value = -15.2778
    $ws.Cells.Item(22, 10).Value = -0.11
    $ws.Cells.Item(22, 11).Value = 99.81999999999999
    $ws.Cells.Item(22, 12).Value = 0
    $ws.Cells.Item(22, 13).Value = 0
    $ws.Cells.Item(22, 14).Value = 0.6
    $ws.Cells.Item(22, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(22, 16).Value = "early_exit"
    $ws.Cells.Item(22, 17).Value = 0.11
}

Add-TradeRow($wb.Worksheets.Item("All Trades"))
Add-TradeRow($wb.Worksheets.Item("MarketMaking"))
